# Move logic from node example to TypeScript SDK tests
# - Delete the extra worksheets "Doesitwork " and "Helloworld!"
# - Update A4/B4 on Sheet1: replace the boolean literal with a formula that
#   checks whether the square (B2) is >= 17
# - Update the selection to A5

$wb = $excel.ActiveWorkbook

# Remove the extra sheets, keeping only Sheet1
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Doesitwork ").Delete()
$wb.Worksheets.Item("Helloworld!").Delete()
$excel.DisplayAlerts = $true

$ws = $wb.Worksheets.Item("Sheet1")

# Update label and formula in row 4
$ws.Range("A4").Value = "Square >= 17?"
$ws.Range("B4").Formula = "=B2>=17"

# Update the active selection
$ws.Range("A5").Select()
